# Fix the "pcs_frequency" -> "psc_frequency" typo in the pscStats template
# sheet, and leave the selection on the corrected cell (A7), matching the
# author's "fixing the same typo in the template sheet..." commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "psc_frequency"
$ws.Range("A7").Select()
